# Adds new "failure screenshot" related test-data rows:
#  - Sheet1: three new data rows (38-40) for mobile numbers 7980000066-68
#  - Test Data: marks the corresponding mobile numbers (rows 67-69) as "used"

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Test Data")

function Set-TextCell($cell, [string]$text) {
    # Force the cell to store a text (shared-string) value even when the
    # text looks numeric, without leaving a lingering custom number format
    # behind on the cell (keeps it on the workbook's default style).
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

# ---- Sheet1: new rows 38, 39, 40 -----------------------------------------
$rows = @(
    @("7980000066", "test151776", "test151776@gmail.com", "SoftSuave52224"),
    @("7980000067", "test151776", "test151776@gmail.com", "SoftSuave52224"),
    @("7980000068", "test151776", "test151776@gmail.com", "SoftSuave52224")
)

$r = 38
foreach ($row in $rows) {
    Set-TextCell $ws1.Cells.Item($r, 1) $row[0]
    Set-TextCell $ws1.Cells.Item($r, 2) $row[1]
    Set-TextCell $ws1.Cells.Item($r, 3) $row[2]
    Set-TextCell $ws1.Cells.Item($r, 4) $row[3]
    $r = $r + 1
}

# ---- Test Data: mark 7980000066/67/68 as "used" ---------------------------
$ws2.Cells.Item(67, 2).Value = "used"
$ws2.Cells.Item(68, 2).Value = "used"
$ws2.Cells.Item(69, 2).Value = "used"
